# Sync attendance_reports: normalize "System"/"system" token position
# in column G ("Recorded By") so that the System/system entry comes first.
#
# Rule observed from the source diff:
#   - Cells with a single value are left untouched.
#   - Cells whose LAST comma-separated token is "System"/"system" (and whose
#     first token is not) have that token moved to the front, e.g.
#       "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   - Cells where BOTH the first and last tokens are "System"/"system"
#     (differing only by case) simply swap those two tokens, e.g.
#       "system, backup@backdoor.com, System"
#         -> "System, backup@backdoor.com, system"
#   - Any other combination (no System/system token in first/last position,
#     e.g. "admin@admin.com, dnasr281@gmail.com") is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $original = $cell.Value2

    if ($null -eq $original -or $original -eq "") {
        continue
    }

    $parts = $original.ToString().Split(",")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Count -le 1) {
        continue
    }

    $lastIdx = $parts.Count - 1
    $first = $parts[0]
    $last = $parts[$lastIdx]
    $firstIsSystem = ($first.ToLower() -eq "system")
    $lastIsSystem = ($last.ToLower() -eq "system")

    $newValue = $null

    if ($lastIsSystem -and -not $firstIsSystem) {
        # Move the trailing System/system token to the front, shifting the
        # rest down by one.
        $reordered = $original.ToString().Split(",")
        for ($i = 0; $i -lt $reordered.Count; $i++) {
            $reordered[$i] = $reordered[$i].Trim()
        }
        for ($i = $lastIdx; $i -gt 0; $i--) {
            $reordered[$i] = $reordered[$i - 1]
        }
        $reordered[0] = $last
        $newValue = [string]::Join(", ", $reordered)
    }
    elseif ($firstIsSystem -and $lastIsSystem) {
        # Both ends are System/system (possibly different case) - swap them.
        $swapped = $original.ToString().Split(",")
        for ($i = 0; $i -lt $swapped.Count; $i++) {
            $swapped[$i] = $swapped[$i].Trim()
        }
        $swapped[0] = $last
        $swapped[$lastIdx] = $first
        $newValue = [string]::Join(", ", $swapped)
    }

    if ($null -ne $newValue -and -not $newValue.Equals($original.ToString())) {
        $cell.Value = $newValue
    }
}
